# Apply the "Figure 4 + associated uncertainty data" update described in the commit:
#   - relabel the "Type" column entries used on the Fig4C / Fig4D uncertainty sheets
#     ("Direct" -> "Land CDR (Direct)", "Indirect" -> "Indirect Removals",
#      "Total" -> "Total Land Removals")
#   - make Fig4C the active/selected sheet (was Fig4D) with C2:C10 selected
#   - mirror the same C2:C10 selection on Fig4D

$wb = $excel.ActiveWorkbook

$newDirect   = "Land CDR (Direct)"
$newIndirect = "Indirect Removals"
$newTotal    = "Total Land Removals"

foreach ($sheetName in @("Fig4C", "Fig4D")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("C2").Value = $newDirect
    $ws.Range("C3").Value = $newIndirect
    $ws.Range("C4").Value = $newTotal

    $ws.Range("C5").Value = $newDirect
    $ws.Range("C6").Value = $newIndirect
    $ws.Range("C7").Value = $newTotal

    $ws.Range("C8").Value = $newDirect
    $ws.Range("C9").Value = $newIndirect
    $ws.Range("C10").Value = $newTotal
}

# Fig4C becomes the active sheet/tab (previously Fig4D was active)
$wsC = $wb.Worksheets.Item("Fig4C")
[void]$wsC.Activate()
[void]$wsC.Range("C2:C10").Select()

# Fig4D keeps the same selected range, but is no longer the active tab
$wsD = $wb.Worksheets.Item("Fig4D")
[void]$wsD.Range("C2:C10").Select()

# Re-activate Fig4C last so it remains the tab shown/selected on open
[void]$wsC.Activate()
